$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (bold, border, center/top alignment) from the existing
# numbered-row label in column A down into the four new rows.
$ws.Range("A11").Copy()
$ws.Range("A12:A15").PasteSpecial(-4122)

# Row 12
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = "HOU"
$ws.Cells.Item(12,3).Value = "GSW"
$ws.Cells.Item(12,4).Value = "away"
$ws.Cells.Item(12,5).Value = "'2025-05-02"
$ws.Cells.Item(12,5).Style = $ws.Cells.Item(2,5).Style
$ws.Cells.Item(12,6).Value = "240:00"
$ws.Cells.Item(12,7).Value = 35
$ws.Cells.Item(12,8).Value = 78
$ws.Cells.Item(12,9).Value = 0.449
$ws.Cells.Item(12,10).Value = 12
$ws.Cells.Item(12,11).Value = 30
$ws.Cells.Item(12,12).Value = 0.4
$ws.Cells.Item(12,13).Value = 33
$ws.Cells.Item(12,14).Value = 46
$ws.Cells.Item(12,15).Value = 0.717
$ws.Cells.Item(12,16).Value = 11
$ws.Cells.Item(12,17).Value = 35
$ws.Cells.Item(12,18).Value = 46
$ws.Cells.Item(12,19).Value = 23
$ws.Cells.Item(12,20).Value = 12
$ws.Cells.Item(12,21).Value = 5
$ws.Cells.Item(12,22).Value = 11
$ws.Cells.Item(12,23).Value = 18
$ws.Cells.Item(12,24).Value = 115
$ws.Cells.Item(12,25).Value = 8
$ws.Cells.Item(12,26).Value = 25
$ws.Cells.Item(12,27).Value = 28
$ws.Cells.Item(12,28).Value = 33
$ws.Cells.Item(12,29).Value = 29
$ws.Cells.Item(12,30).Value = "W"

# Row 13
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = "GSW"
$ws.Cells.Item(13,3).Value = "HOU"
$ws.Cells.Item(13,4).Value = "home"
$ws.Cells.Item(13,5).Value = "'2025-05-02"
$ws.Cells.Item(13,5).Style = $ws.Cells.Item(2,5).Style
$ws.Cells.Item(13,6).Value = "240:00"
$ws.Cells.Item(13,7).Value = 37
$ws.Cells.Item(13,8).Value = 90
$ws.Cells.Item(13,9).Value = 0.411
$ws.Cells.Item(13,10).Value = 15
$ws.Cells.Item(13,11).Value = 49
$ws.Cells.Item(13,12).Value = 0.306
$ws.Cells.Item(13,13).Value = 18
$ws.Cells.Item(13,14).Value = 22
$ws.Cells.Item(13,15).Value = 0.818
$ws.Cells.Item(13,16).Value = 13
$ws.Cells.Item(13,17).Value = 29
$ws.Cells.Item(13,18).Value = 42
$ws.Cells.Item(13,19).Value = 26
$ws.Cells.Item(13,20).Value = 6
$ws.Cells.Item(13,21).Value = 6
$ws.Cells.Item(13,22).Value = 16
$ws.Cells.Item(13,23).Value = 30
$ws.Cells.Item(13,24).Value = 107
$ws.Cells.Item(13,25).Value = -8
$ws.Cells.Item(13,26).Value = 21
$ws.Cells.Item(13,27).Value = 27
$ws.Cells.Item(13,28).Value = 36
$ws.Cells.Item(13,29).Value = 23
$ws.Cells.Item(13,30).Value = "L"

# Row 14
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = "GSW"
$ws.Cells.Item(14,3).Value = "HOU"
$ws.Cells.Item(14,4).Value = "away"
$ws.Cells.Item(14,5).Value = "'2025-05-04"
$ws.Cells.Item(14,5).Style = $ws.Cells.Item(2,5).Style
$ws.Cells.Item(14,6).Value = "240:00"
$ws.Cells.Item(14,7).Value = 39
$ws.Cells.Item(14,8).Value = 82
$ws.Cells.Item(14,9).Value = 0.476
$ws.Cells.Item(14,10).Value = 18
$ws.Cells.Item(14,11).Value = 43
$ws.Cells.Item(14,12).Value = 0.419
$ws.Cells.Item(14,13).Value = 7
$ws.Cells.Item(14,14).Value = 9
$ws.Cells.Item(14,15).Value = 0.778
$ws.Cells.Item(14,16).Value = 3
$ws.Cells.Item(14,17).Value = 35
$ws.Cells.Item(14,18).Value = 38
$ws.Cells.Item(14,19).Value = 24
$ws.Cells.Item(14,20).Value = 6
$ws.Cells.Item(14,21).Value = 8
$ws.Cells.Item(14,22).Value = 7
$ws.Cells.Item(14,23).Value = 14
$ws.Cells.Item(14,24).Value = 103
$ws.Cells.Item(14,25).Value = 14
$ws.Cells.Item(14,26).Value = 23
$ws.Cells.Item(14,27).Value = 28
$ws.Cells.Item(14,28).Value = 19
$ws.Cells.Item(14,29).Value = 33
$ws.Cells.Item(14,30).Value = "W"

# Row 15
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = "HOU"
$ws.Cells.Item(15,3).Value = "GSW"
$ws.Cells.Item(15,4).Value = "home"
$ws.Cells.Item(15,5).Value = "'2025-05-04"
$ws.Cells.Item(15,5).Style = $ws.Cells.Item(2,5).Style
$ws.Cells.Item(15,6).Value = "240:00"
$ws.Cells.Item(15,7).Value = 34
$ws.Cells.Item(15,8).Value = 84
$ws.Cells.Item(15,9).Value = 0.405
$ws.Cells.Item(15,10).Value = 6
$ws.Cells.Item(15,11).Value = 18
$ws.Cells.Item(15,12).Value = 0.333
$ws.Cells.Item(15,13).Value = 15
$ws.Cells.Item(15,14).Value = 21
$ws.Cells.Item(15,15).Value = 0.714
$ws.Cells.Item(15,16).Value = 14
$ws.Cells.Item(15,17).Value = 38
$ws.Cells.Item(15,18).Value = 52
$ws.Cells.Item(15,19).Value = 14
$ws.Cells.Item(15,20).Value = 5
$ws.Cells.Item(15,21).Value = 5
$ws.Cells.Item(15,22).Value = 11
$ws.Cells.Item(15,23).Value = 14
$ws.Cells.Item(15,24).Value = 89
$ws.Cells.Item(15,25).Value = -14
$ws.Cells.Item(15,26).Value = 19
$ws.Cells.Item(15,27).Value = 20
$ws.Cells.Item(15,28).Value = 23
$ws.Cells.Item(15,29).Value = 27
$ws.Cells.Item(15,30).Value = "L"
